$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 173 -----------------------------------------------------------
$ws.Cells.Item(173, 2).Value = 0
$ws.Cells.Item(173, 3).Value = 2.58999991416931
$ws.Cells.Item(173, 4).Value = 2.58999991416931
$ws.Cells.Item(173, 5).Value = 2.58999991416931
$ws.Cells.Item(173, 6).Value = 2.58999991416931

# A173: same date style as A172 (numFmt "yyyy-mm-dd hh:mm:ss"), value 2024-06-04 07:00
$ws.Cells.Item(173, 1).Value = 45447.2916666667
$ws.Cells.Item(172, 1).Copy()
$ws.Cells.Item(173, 1).PasteSpecial(-4122)

# G173: reuse the existing text value "2.58999991416931" (same as G172, shared string)
$ws.Cells.Item(172, 7).Copy()
$ws.Cells.Item(173, 7).PasteSpecial(-4104)

# H173: ticker text "EAV.MI" (same as H172, shared string)
$ws.Cells.Item(172, 8).Copy()
$ws.Cells.Item(173, 8).PasteSpecial(-4104)

# --- Row 174 -----------------------------------------------------------
$ws.Cells.Item(174, 2).Value = 66000
$ws.Cells.Item(174, 3).Value = 3.09999990463257
$ws.Cells.Item(174, 4).Value = 2.63000011444092
$ws.Cells.Item(174, 5).Value = 2.63000011444092
$ws.Cells.Item(174, 6).Value = 2.97000002861023

# A174: same date style as A172, value 2024-06-05 15:28:48
$ws.Cells.Item(174, 1).Value = 45448.645
$ws.Cells.Item(172, 1).Copy()
$ws.Cells.Item(174, 1).PasteSpecial(-4122)

# G174: reuse the existing text value "2.97000002861023" (already present as shared string)
$ws.Cells.Item(27, 7).Copy()
$ws.Cells.Item(174, 7).PasteSpecial(-4104)

# H174: ticker text "EAV.MI"
$ws.Cells.Item(172, 8).Copy()
$ws.Cells.Item(174, 8).PasteSpecial(-4104)

$excel.CutCopyMode = $false
